$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(622, "Milwaukee Bucks", 141, "Detroit Pistons", 135, "No", 17832, "Little Caesars Arena", "Milwaukee Bucks", "Detroit Pistons"),
    @(623, "Philadelphia 76ers", 97, "Charlotte Hornets", 89, "No", 17832, "Spectrum Center", "Philadelphia 76ers", "Charlotte Hornets"),
    @(624, "San Antonio Spurs", 131, "Washington Wizards", 127, "No", 17832, "Capital One Arena", "San Antonio Spurs", "Washington Wizards"),
    @(625, "Cleveland Cavaliers", 116, "Atlanta Hawks", 95, "No", 17832, "State Farm Arena", "Cleveland Cavaliers", "Atlanta Hawks"),
    @(626, "Toronto Raptors", 100, "New York Knicks", 126, "No", 17832, "Madison Square Garden (IV)", "New York Knicks", "Toronto Raptors"),
    @(627, "Memphis Grizzlies", 96, "Chicago Bulls", 125, "No", 17832, "United Center", "Chicago Bulls", "Memphis Grizzlies"),
    @(628, "Utah Jazz", 126, "Houston Rockets", 127, "No", 17832, "Toyota Center", "Houston Rockets", "Utah Jazz"),
    @(629, "Oklahoma City Thunder", 102, "Minnesota Timberwolves", 97, "No", 17832, "Target Center", "Oklahoma City Thunder", "Minnesota Timberwolves")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 2).NumberFormat = "#,##0"
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 4).NumberFormat = "#,##0"
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
}

"Done writing rows 622-629"
